$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several updated "Price" values are plain decimals (e.g. "227.16") that Excel
# would otherwise auto-convert to numbers. The sheet stores Price as text, so
# mark those specific cells as Text before writing the new value, preserving
# their original text representation (leading zeros, trailing zeros, etc).
$textForceRows = @(5, 6, 7, 8, 10, 12, 14, 17, 18, 20, 21, 22, 23, 24, 25, 26, 27, 28, 31, 32, 33, 35, 36, 37, 38, 40, 41, 43, 44, 45, 46, 47, 48, 49, 50, 51)
foreach ($r in $textForceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "92.464.33"
$ws.Range("E2").Value = "  -2.89%  "

$ws.Range("D3").Value = "3.291.17"
$ws.Range("E3").Value = "  -4.71%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "227.16"
$ws.Range("E5").Value = "  -5.33%  "

$ws.Range("D6").Value = "606.78"
$ws.Range("E6").Value = "  -5.67%  "

$ws.Range("D7").Value = "1.35"
$ws.Range("E7").Value = "  -8.16%  "

$ws.Range("D8").Value = "0.376"
$ws.Range("E8").Value = "  -6.67%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "0.913"
$ws.Range("E10").Value = "  -9.02%  "

$ws.Range("D11").Value = "3.285.01"
$ws.Range("E11").Value = "  -4.90%  "

$ws.Range("D12").Value = "41.50"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("E13").Value = "  -3.43%  "

$ws.Range("D14").Value = "5.88"
$ws.Range("E14").Value = "  -3.71%  "

$ws.Range("D15").Value = "92.253.44"
$ws.Range("E15").Value = "  -2.80%  "

$ws.Range("D16").Value = "3.901.00"
$ws.Range("E16").Value = "  -4.88%  "

$ws.Range("D17").Value = "0.0000239"
$ws.Range("E17").Value = "  -6.92%  "

$ws.Range("D18").Value = "7.89"
$ws.Range("E18").Value = "  -6.43%  "

$ws.Range("D19").Value = "3.279.84"
$ws.Range("E19").Value = "  -4.57%  "

$ws.Range("D20").Value = "16.95"
$ws.Range("E20").Value = "  -4.84%  "

$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  -7.23%  "

$ws.Range("D22").Value = "3.38"
$ws.Range("E22").Value = "  +6.43%  "

$ws.Range("D23").Value = "482.56"
$ws.Range("E23").Value = "  -3.91%  "

$ws.Range("D24").Value = "0.436"
$ws.Range("E24").Value = "  -13.63%  "

$ws.Range("D25").Value = "0.0000175"
$ws.Range("E25").Value = "  -8.86%  "

$ws.Range("D26").Value = "5.97"
$ws.Range("E26").Value = "  -9.36%  "

$ws.Range("D27").Value = "88.65"
$ws.Range("E27").Value = "  -3.48%  "

$ws.Range("D28").Value = "11.53"
$ws.Range("E28").Value = "  -4.07%  "

$ws.Range("D29").Value = "3.468.01"
$ws.Range("E29").Value = "  -4.56%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").Value = "10.86"
$ws.Range("E31").Value = "  -7.50%  "

$ws.Range("D32").Value = "0.135"
$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("D33").Value = "2.57"
$ws.Range("E33").Value = "  -6.50%  "

$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("D35").Value = "0.170"
$ws.Range("E35").Value = "  -7.67%  "

$ws.Range("D36").Value = "27.71"
$ws.Range("E36").Value = "  -10.50%  "

$ws.Range("D37").Value = "0.517"
$ws.Range("E37").Value = "  -8.32%  "

$ws.Range("D38").Value = "535.45"
$ws.Range("E38").Value = "  +1.92%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "7.19"
$ws.Range("E40").Value = "  -6.50%  "

$ws.Range("D41").Value = "0.145"
$ws.Range("E41").Value = "  -3.67%  "

$ws.Range("E42").Value = "  -7.26%  "

$ws.Range("D47").Value = "0.0402"
$ws.Range("E47").Value = "  -3.19%  "

$ws.Range("D48").Value = "5.20"
$ws.Range("E48").Value = "  -7.79%  "

$ws.Range("D51").Value = "7.72"
$ws.Range("E51").Value = "  -3.84%  "

# Reordered rows (coin pairs swapped between adjacent ranking rows)
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "23.87"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "0.843"
$ws.Range("E44").Value = "  -7.89%  "

$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "3.55"
$ws.Range("E45").Value = "  +1.65%  "

$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "1.64"
$ws.Range("E46").Value = "  -3.27%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "51.47"
$ws.Range("E49").Value = "  -3.81%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "2.06"
$ws.Range("E50").Value = "  -4.20%  "

